# edit.ps1 - Apply updated vm_pu values (case with 380 kV) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$colLetters = @('B','C','D','E','F','I','J','K','L','M','N')

$data = @(
    ,@(2, 1.02, 1.046498230165854, 1.05209406302, 1.044061878458996, 1.059866396611977, 1.043381393377614, 1.051552304943406, 1.054843747740171, 1.046833987233666, 1.062594733175358, 1.053045629831161)
    ,@(3, 1.02, 1.047825331144027, 1.052889814815167, 1.045203821032957, 1.060927931873474, 1.043714129534724, 1.052525638492577, 1.05545232037632, 1.047786227858745, 1.063469956949974, 1.054020345625614)
    ,@(4, 1.02, 1.048683522913168, 1.053404517575569, 1.045942514299506, 1.06161491885446, 1.043928042891626, 1.053154450728691, 1.055845250439352, 1.048401592720704, 1.064035786435757, 1.05465005084723)
    ,@(5, 1.02, 1.049044182985811, 1.053620850568019, 1.04625301044205, 1.061903753724952, 1.044017640075369, 1.053418566227334, 1.056010233562708, 1.048660102921622, 1.064273542676538, 1.054914541420175)
    ,@(6, 1.02, 1.049104732242517, 1.053657171019764, 1.046305141211187, 1.061952251896366, 1.044032664385517, 1.053462898536728, 1.056037922967937, 1.048703496885129, 1.064313456075441, 1.054958936686534)
    ,@(7, 1.02, 1.048688342552655, 1.053407408416065, 1.045946663359947, 1.061618778178576, 1.043929241397357, 1.053157980782134, 1.055847455755513, 1.048405047687585, 1.064038963810641, 1.054653585913753)
    ,@(8, 1.02, 1.046946842242331, 1.052363033031428, 1.044447850132393, 1.060225126742871, 1.04349413111138, 1.051881455857084, 1.055049595502003, 1.047155967935823, 1.062890623001984, 1.053375248176891)
    ,@(9, 1.02, 1.043873872385652, 1.050521150895609, 1.041804968660973, 1.057770077548897, 1.042716745317592, 1.049624302818434, 1.053637075250404, 1.044948725445096, 1.060863236515641, 1.051114889722026)
    ,@(10, 1.02, 1.041822167848874, 1.049292159605881, 1.040041712777995, 1.056133813259762, 1.04219127880502, 1.048114183920264, 1.052690929537056, 1.043472932670255, 1.059508985467867, 1.049602626281808)
    ,@(11, 1.02, 1.040932982696108, 1.048759733345846, 1.039277855585875, 1.055425380758042, 1.041962027129596, 1.047458985513404, 1.052280170159948, 1.042832849968773, 1.058921934121521, 1.048946497418058)
    ,@(12, 1.02, 1.040602577700281, 1.048561926132229, 1.038994068776345, 1.055162248097219, 1.041876613425126, 1.047215416511979, 1.05212743394182, 1.042594933731784, 1.058703777830657, 1.048702582520713)
    ,@(13, 1.02, 1.040673456285533, 1.048604358288578, 1.039054944600119, 1.055218690484245, 1.041894946707612, 1.047267671945964, 1.052160203718719, 1.042645974900163, 1.058750577610451, 1.048754912163409)
    ,@(14, 1.02, 1.040905673830687, 1.048743383367361, 1.039254398843081, 1.055403629919895, 1.041954972103158, 1.047438856102319, 1.052267548244004, 1.042813187038092, 1.058903903287759, 1.0489263394209)
    ,@(15, 1.02, 1.041048734423542, 1.048829035936253, 1.039377281649306, 1.055517578587041, 1.041991921329159, 1.047544301839294, 1.052333665270907, 1.042916190541015, 1.058998359100879, 1.049031934902916)
    ,@(16, 1.02, 1.041881163262174, 1.049327489351285, 1.040092399730955, 1.056180831172076, 1.042206457148835, 1.048157639507159, 1.052718167626407, 1.043515390411158, 1.05954793228483, 1.049646143580621)
    ,@(17, 1.02, 1.042403111056593, 1.049640084666657, 1.040540877198316, 1.056596892613509, 1.042340568397461, 1.048542018149526, 1.052959068192978, 1.04389096847954, 1.059892489578163, 1.050031068084749)
    ,@(18, 1.02, 1.042707479066892, 1.049822390706014, 1.04080243238192, 1.056839582100113, 1.042418627205172, 1.048766093795545, 1.053099478039093, 1.044109935225035, 1.060093401231917, 1.050255461943887)
    ,@(19, 1.02, 1.042811248049381, 1.049884548043057, 1.040891610234109, 1.056922334323357, 1.042445215107346, 1.048842476526864, 1.053147336658092, 1.044184580072725, 1.06016189630514, 1.050331953147445)
    ,@(20, 1.02, 1.042347118820019, 1.049606548792273, 1.040492763362534, 1.056552252353813, 1.04232619670785, 1.048500791003918, 1.05293323254165, 1.043850683034514, 1.059855528346467, 1.049989782391865)
    ,@(21, 1.02, 1.040837294955159, 1.048702445066661, 1.039195666114693, 1.055349169581277, 1.041937303288366, 1.0473884521546, 1.052235942432673, 1.042763951696791, 1.058858755449072, 1.048875863893795)
    ,@(22, 1.02, 1.039887300535375, 1.048133765160456, 1.038379802829231, 1.0545928049667, 1.041691289193284, 1.046687927018796, 1.051796591364321, 1.042079747669683, 1.058231469663559, 1.048174343931946)
    ,@(23, 1.02, 1.04039097912947, 1.048435255520002, 1.038812339300795, 1.05499376282654, 1.041821848510881, 1.047059398923338, 1.052029588785029, 1.042442546457743, 1.058564060699914, 1.048546343369201)
    ,@(24, 1.02, 1.042372419534076, 1.049621702289358, 1.040514504039315, 1.056572423341381, 1.042332691168733, 1.048519420165391, 1.052944906879897, 1.043868886606885, 1.059872229730158, 1.050008438008884)
    ,@(25, 1.02, 1.044668831708335, 1.050997508537227, 1.042488442044902, 1.058404685170926, 1.042918986058251, 1.050208763877276, 1.054003030027271, 1.045520098755983, 1.061387828749367, 1.051700180782611)
)

foreach ($rowData in $data) {
    $rowNum = $rowData[0]
    for ($i = 0; $i -lt $colLetters.Count; $i++) {
        $col = $colLetters[$i]
        $val = $rowData[$i + 1]
        $ws.Range("$col$rowNum").Value2 = $val
    }
}
